# edit.ps1 - Applies the diff described in the task:
#  1. Rewrites run/proofErr structure inside the "Latency window" paragraph.
#  2. Turns an empty paragraph after "...type 5." into one carrying an
#     eastAsia rFonts hint (pPr/rPr only, still empty text).
#  3. Merges the "// Check available Nack..." paragraph and the following
#     empty paragraph into a single empty (but eastAsia-hinted) paragraph.
#  4. Removes a stray <w:lastRenderedPageBreak/> from the run that starts
#     "Otherwise keep into the buffer".
#  5. Inserts a brand-new paragraph "LatencyWindow " (with "Lat" split into
#     its own eastAsia-hinted run) right after the ACKACK_TS paragraph.
#
# NOTE: Word's COM object model renumbers Paragraphs(...) live, so hunks
# are applied top-to-bottom and later indices account for paragraph-count
# deltas introduced by earlier hunks (hunk 3 removes one paragraph, so the
# hunk 4 / hunk 5 target indices below are their *original* index minus 1).

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns = "http://schemas.microsoft.com/office/word/2010/wordml"

# ---------------------------------------------------------------------
# Hunk 1 (original paragraph 3): restructure the "Latency window" runs.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(3)
$full1 = $p1.Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6EA1E91F" w14:textId="4ACEEE43" w:rsidR="009F2F30" w:rsidRDefault="009F2F30"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t>In the request packet, receiver will includes Window size</w:t></w:r><w:r w:rsidR="00D55691"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D55691"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>and</w:t></w:r><w:r><w:t xml:space="preserve"> Latency window </w:t></w:r><w:r><w:t xml:space="preserve">to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>help  sender</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> initialize the window.</w:t></w:r></w:p>
'@
$full1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Hunk 2 (original paragraph 17): give the empty paragraph after
# "...type 5." a pPr/rPr eastAsia rFonts hint.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(17)
$full2 = $p2.Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="53428FC5" w14:textId="77777777" w:rsidR="0059720F" w:rsidRDefault="0059720F"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>
'@
$full2.InsertXML($xml2)

# ---------------------------------------------------------------------
# Hunk 3 (original paragraphs 32-33): delete the "// Check available
# Nack..." text (and the empty paragraph after it), leaving a single
# empty paragraph with an eastAsia rFonts hint.
# ---------------------------------------------------------------------
$p3a = $d.Paragraphs(32)
$p3b = $d.Paragraphs(33)
$full3 = $d.Range($p3a.Range.Start, $p3b.Range.End)
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="55FA4F21" w14:textId="4971615D" w:rsidR="00754153" w:rsidRDefault="00754153"><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>
'@
$full3.InsertXML($xml3)

# ---------------------------------------------------------------------
# Hunk 4 (original paragraph 37, now 36 because hunk 3 merged two
# paragraphs into one): strip <w:lastRenderedPageBreak/> from the run
# that reads "Otherwise keep into the buffer".
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(36)
$full4 = $p4.Range
$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5E58A5D1" w14:textId="268E851C" w:rsidR="00754153" w:rsidRDefault="00754153"><w:r><w:t>Otherwise keep into the buffer</w:t></w:r></w:p>
'@
$full4.InsertXML($xml4)

# ---------------------------------------------------------------------
# Hunk 5 (original paragraph 73, now 72): insert a new paragraph
# "LatencyWindow " right after the ACKACK_TS paragraph, before the
# following empty paragraph.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(72)
$full5 = $p5.Range
$collapsed5 = $d.Range($full5.End - 1, $full5.End - 1)
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Lat</w:t></w:r><w:r><w:t>encyWindow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$collapsed5.InsertXML($xml5)

Write-Output ("Done. Final paragraph count=" + $d.Paragraphs.Count)
